$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-11 14:45:23"
$wsZh.Range("H2").Value = "2016-03-11 14:45:41"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-11 14:45:26"
$wsDe.Range("H2").Value = "2016-03-11 14:45:50"
